$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("SeokkSquirrel", "afk"),
    @("WiX21", "afk"),
    @("Monkeygoberserk", "afk"),
    @("Khabib_Time", "ak"),
    @("Enaxie", "afk"),
    @("BlizzlerButNot", "afk"),
    @("GRANDMAA", "afk"),
    @("gosling123", "afk"),
    @("UNDOMINABLE", "throwing game"),
    @("imasavage", "afk"),
    @("TrashOfCountsFamily", "afk"),
    @("Mkool14", "afk"),
    @("Polaris0", "afk"),
    @("z4ra0", "afk"),
    @("Sephiroth_99", "afk"),
    @("khaiworld", "afk"),
    @("chug", "afk"),
    @("ice_dragon95", "homophobia"),
    @("BigPapaPanther", "afk"),
    @("N3xus", "afk"),
    @("Ekizeel", "afk")
)

$row = 6
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $row++
}

# Extend the used range down to row 27 (matching the original trailing blank
# row pattern) without introducing any new cell style. Touching the border
# property with "no line" dirties the cell for used-range tracking purposes
# but resolves back to the default style, so no new cellXfs entry is created.
$ws.Range("A27").Borders.LineStyle = -4142

$null = $ws.Range("B10").Select()

